# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Replace the old "Estado de Cuenta" worker/period detail rows with a new
# data set, update the summary totals, and drop the now-unused trailing rows
# (the footer block moves up to sit right under the last data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last populated data row (old row 35) carries the "closing" border
# formatting that the new last data row (row 29) must inherit. Copy that
# formatting over before we touch anything else.
$ws.Range("B35:J35").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)

# New worker/period dataset (replaces rows 16-35 worth of old data with a
# fresh 14-row table in rows 16-29).
$data = @(
  @("CC","73103129","WILSON EDUARDO TORREGLOSA PAUT","2411",52000,1300000),
  @("CC","73201649","JULIO CESAR ALMAGRO PALENCIA","2411",52000,1300000),
  @("CC","1103103025","MARIA MARGARITA HERNANDEZ BADEL","2411",160000,4000000),
  @("CC","1044919362","RICHARD LEONARDO FERNANDEZ YEPEZ","2411",52000,1300000),
  @("CC","1143334657","HENRY IRIARTE ORTEGA","2411",52000,1300000),
  @("CC","1128049624","JWAUIS BELTRAN PENAGOS","2411",52000,1056200),
  @("CC","73207278","FRANCISCO JAVIER BARBOZA ORTEGA","2411",52000,1300000),
  @("CC","1002315865","EDUARD DAVID VILLADIEGO MORALES","2411",52000,1300000),
  @("CC","1143390615","BRAYAN IRIARTE ALMAGRO","2411",52000,1300000),
  @("CC","91324881","JUAN PABLO GARCIA PINTO","2403",52000,1300000),
  @("CC","91324881","JUAN PABLO GARCIA PINTO","2402",52000,1300000),
  @("CC","1033791107","JEISSON ORLANDO RIAÃ?O LEON","2407",140000,3200000),
  @("CC","1033791107","JEISSON ORLANDO RIAÃ?O LEON","2406",140000,3200000),
  @("CC","1033791107","JEISSON ORLANDO RIAÃ?O LEON","2405",18667,3200000)
)

$r = 16
foreach ($row in $data) {
  $ws.Cells.Item($r,2).Value = $row[0]
  $ws.Cells.Item($r,3).Value = $row[1]
  $ws.Cells.Item($r,4).Value = $row[2]
  $ws.Cells.Item($r,5).Value = $row[3]
  $ws.Cells.Item($r,6).Value = $row[4]
  $ws.Cells.Item($r,7).Value = $row[5]
  $r = $r + 1
}

# Totals: sum of "Valor Mora" for the new data set, and the new distinct
# period count (11 distinct workers is unchanged; periods drop from 12 to 6).
$ws.Range("E11").Value = 978667
$ws.Range("F13").Value = 6

# Drop the now-empty/old rows 30-35 so the signature-block footer (old rows
# 40-41) slides up to sit directly under the new data (new rows 34-35).
$ws.Rows("30:35").Delete()
